# LOB1254.xlsx edit: add the "Objectives / Docentes responsáveis / Programa
# resumido / Programa / Método / Critério / Norma de recuperação /
# Bibliografia" section content that was missing, shifting the lower half
# of the table down by one row to make room for the teacher's name on its
# own row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room: insert a new row at 13 (everything from the old row 13
#    down shifts to 14.. and keeps its formatting/row-height).
$ws.Rows(13).Insert()

# Copy the wrap-text formatting from the (now shifted) B14:C14 cells into
# the freshly inserted B13:C13 so they pick up styles 2/3 like every other
# data row, then clear out the leftover label formatting in A13 (this row
# has no label in column A).
$ws.Range("B14:C14").Copy() | Out-Null
$ws.Range("B13:C13").PasteSpecial(-4122) | Out-Null
$ws.Range("A13").Clear()

# 2) Fill in the actual values.

# Objetivos: (row 10) previously held the teacher's name by mistake; it
# should hold the real Portuguese objective text.
$ws.Range("B10").Value = "Propiciar conhecimentos básicos sobre os materiais terrestres e os principais processos geológicos."
$ws.Range("C10").Value = "Propiciar conhecimentos básicos sobre os materiais terrestres e os principais processos geológicos."

# Docentes responsáveis: (row 12 label) -> new row 13 holds the teacher.
$ws.Range("B13").Value = "5464150 - Mariana Consiglio Kasemodel"
$ws.Range("C13").Value = "5464150 - Mariana Consiglio Kasemodel"

# Programa resumido: (row 14) previously held "Semestral" by mistake.
$ws.Range("B14").Value = "Processos endógenos e exógenos da Terra. Materiais constituintes da crosta terrestre (minerais e rochas)."
$ws.Range("C14").Value = "Processos endógenos e exógenos da Terra. Materiais constituintes da crosta terrestre (minerais e rochas)."

# Programa: (row 16) previously held a stray date value.
$ws.Range("B16").Value = "Breve história da Geologia. Materiais constituintes da crosta terrestre (minerais e rochas). Origem e constituição do universo, do sistema solar e da Terra. Estrutura interna da Terra. Composição da Terra. Processos endógenos e exógenos (dinâmica interna e externa da Terra).  Teoria da tectônica de placas.  Rochas ígneas e vulcanismo. Rochas metamórficas e metamorfismo. Rochas sedimentares. Intemperismo, erosão, transporte de sedimentos.  Estrutura geológicas. Tempo geológico e estratigrafia."
$ws.Range("C16").Value = "Breve história da Geologia. Materiais constituintes da crosta terrestre (minerais e rochas). Origem e constituição do universo, do sistema solar e da Terra. Estrutura interna da Terra. Composição da Terra. Processos endógenos e exógenos (dinâmica interna e externa da Terra).  Teoria da tectônica de placas.  Rochas ígneas e vulcanismo. Rochas metamórficas e metamorfismo. Rochas sedimentares. Intemperismo, erosão, transporte de sedimentos.  Estrutura geológicas. Tempo geológico e estratigrafia."

# Método: (row 19) previously held the teacher's name by mistake.
$ws.Range("B19").Value = "Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas."
$ws.Range("C19").Value = "Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas."

# Critério: (row 20) previously held the "Método" text.
$ws.Range("B20").Value = "Média ponderada de provas  e atividades."
$ws.Range("C20").Value = "Média ponderada de provas  e atividades."

# Norma de recuperação: (row 21) previously held the "Critério" text.
$ws.Range("B21").Value = "1 (uma) prova escrita"
$ws.Range("C21").Value = "1 (uma) prova escrita"

# Bibliografia: (row 22) previously held the "Norma de recuperação" text.
$ws.Range("B22").Value = "Bibliografia básica:PRESS, F.; SIEVER, R.; GROTZINGER, J.; JORDAN, T. H. Para entender a Terra. Porto Alegre: Bookman, 2008. 656p.REED, W.; MONROE, J. S. Fundamentos de Geologia. São Paulo: Cengage Learning, 2011. 508p.Bibliografia complementar:TEIXEIRA, W.; FAIRCHILD, T. R.; DE TOLEDO, M. C. M.; TAIOLI, F. Decifrando a Terra. São Paulo: Companhia Editora Nacional, 2003. 623p."
$ws.Range("C22").Value = "Bibliografia básica:PRESS, F.; SIEVER, R.; GROTZINGER, J.; JORDAN, T. H. Para entender a Terra. Porto Alegre: Bookman, 2008. 656p.REED, W.; MONROE, J. S. Fundamentos de Geologia. São Paulo: Cengage Learning, 2011. 508p.Bibliografia complementar:TEIXEIRA, W.; FAIRCHILD, T. R.; DE TOLEDO, M. C. M.; TAIOLI, F. Decifrando a Terra. São Paulo: Companhia Editora Nacional, 2003. 623p."

Write-Output "Edit complete"
